$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Value = "58.580.49" },
    @{ Addr = "E2"; Value = "  -0.99%  " },
    @{ Addr = "D3"; Value = "2.488.20" },
    @{ Addr = "E3"; Value = "  -1.44%  " },
    @{ Addr = "E4"; Value = "  +0.16%  " },
    @{ Addr = "D5"; Value = "527.47" },
    @{ Addr = "E5"; Value = "  -1.86%  " },
    @{ Addr = "D6"; Value = "134.27" },
    @{ Addr = "E6"; Value = "  -2.99%  " },
    @{ Addr = "E7"; Value = "  +0.31%  " },
    @{ Addr = "E8"; Value = "  -0.82%  " },
    @{ Addr = "D9"; Value = "0.101" },
    @{ Addr = "E9"; Value = "  -0.59%  " },
    @{ Addr = "E10"; Value = "  -2.07%  " },
    @{ Addr = "D11"; Value = "5.41" },
    @{ Addr = "E11"; Value = "  +0.79%  " },
    @{ Addr = "D12"; Value = "0.343" },
    @{ Addr = "E12"; Value = "  -1.03%  " },
    @{ Addr = "D13"; Value = "2.929.16" },
    @{ Addr = "E13"; Value = "  -1.10%  " },
    @{ Addr = "D14"; Value = "58.464.98" },
    @{ Addr = "E14"; Value = "  -0.94%  " },
    @{ Addr = "D15"; Value = "22.52" },
    @{ Addr = "E15"; Value = "  -2.98%  " },
    @{ Addr = "E16"; Value = "  -1.87%  " },
    @{ Addr = "D17"; Value = "2.490.54" },
    @{ Addr = "E17"; Value = "  -1.45%  " },
    @{ Addr = "D18"; Value = "10.96" },
    @{ Addr = "E18"; Value = "  -1.25%  " },
    @{ Addr = "D19"; Value = "4.22" },
    @{ Addr = "E19"; Value = "  -1.81%  " },
    @{ Addr = "D20"; Value = "322.35" },
    @{ Addr = "E20"; Value = "  -1.03%  " },
    @{ Addr = "E21"; Value = "  +0.05%  " },
    @{ Addr = "D22"; Value = "5.84" },
    @{ Addr = "E22"; Value = "  -0.99%  " },
    @{ Addr = "D23"; Value = "64.49" },
    @{ Addr = "E23"; Value = "  -1.82%  " },
    @{ Addr = "E24"; Value = "  -2.34%  " },
    @{ Addr = "E25"; Value = "  -2.21%  " },
    @{ Addr = "E26"; Value = "  -0.12%  " },
    @{ Addr = "D27"; Value = "7.47" },
    @{ Addr = "E27"; Value = "  -2.49%  " },
    @{ Addr = "D28"; Value = "0.0₃0755" },
    @{ Addr = "E28"; Value = "  -2.95%  " },
    @{ Addr = "D29"; Value = "6.46" },
    @{ Addr = "E29"; Value = "  -4.32%  " },
    @{ Addr = "E30"; Value = "  -3.13%  " },
    @{ Addr = "D31"; Value = "166.93" },
    @{ Addr = "E31"; Value = "  -1.67%  " },
    @{ Addr = "D32"; Value = "1.13" },
    @{ Addr = "E32"; Value = "  -5.52%  " },
    @{ Addr = "E33"; Value = "  +0.04%  " },
    @{ Addr = "D34"; Value = "0.998" },
    @{ Addr = "E34"; Value = "  +0.05%  " },
    @{ Addr = "D35"; Value = "18.26" },
    @{ Addr = "E35"; Value = "  -1.63%  " },
    @{ Addr = "E36"; Value = "  -8.98%  " },
    @{ Addr = "D37"; Value = "4.00" },
    @{ Addr = "E37"; Value = "  -2.93%  " },
    @{ Addr = "D38"; Value = "1.50" },
    @{ Addr = "E38"; Value = "  -4.31%  " },
    @{ Addr = "D39"; Value = "0.803" },
    @{ Addr = "E39"; Value = "  -2.99%  " },
    @{ Addr = "D40"; Value = "3.54" },
    @{ Addr = "E40"; Value = "  -2.86%  " },
    @{ Addr = "E41"; Value = "  -2.83%  " },
    @{ Addr = "D42"; Value = "4.97" },
    @{ Addr = "E42"; Value = "  -5.58%  " },
    @{ Addr = "D43"; Value = "0.598" },
    @{ Addr = "E43"; Value = "  -1.04%  " },
    @{ Addr = "D44"; Value = "127.74" },
    @{ Addr = "E44"; Value = "  -2.42%  " },
    @{ Addr = "D45"; Value = "0.0916" },
    @{ Addr = "E45"; Value = "  -1.83%  " },
    @{ Addr = "D46"; Value = "0.0497" },
    @{ Addr = "E46"; Value = "  -2.91%  " },
    @{ Addr = "D47"; Value = "0.0217" },
    @{ Addr = "E47"; Value = "  -2.25%  " },
    @{ Addr = "D48"; Value = "17.25" },
    @{ Addr = "E48"; Value = "  -1.58%  " },
    @{ Addr = "D49"; Value = "1.741.12" },
    @{ Addr = "E49"; Value = "  -1.37%  " },
    @{ Addr = "D50"; Value = "0.978" },
    @{ Addr = "E50"; Value = "  -0.97%  " },
    @{ Addr = "D51"; Value = "4.70" },
    @{ Addr = "E51"; Value = "  -1.66%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = $origStyle
}